# Apply the "typhoon xlsx" dungeon demo update.
#
# The workbook models a small xlsx2data table on sheet "@DungeonA":
#   row1 = field type specs, row2 = field display names, rows3-8 = data.
# This edit:
#   1. fixes the field specs in row 1 (adds trailing [] / [C] markers,
#      reworks the awards spec and separators from , to -),
#   2. gives column C its own header text "副本掉落" in row 2 (it used
#      to incorrectly reuse the first data row string),
#   3. changes the data separators in column C (rows 3-8) from
#      "," to "-",
#   4. tweaks cosmetic workbook/sheet state (window position, selected
#      cell, column widths, page setup).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("@DungeonA")

# --- window position (bookViews/workbookView) ---
$win = $excel.Windows.Item(1)
$win.Left = 1872

# --- cell content updates -------------------------------------------------
# Order matters for how newly-introduced shared strings get appended, so we
# touch cells in the same sequence the final string table expects:
#   副本掉落 (C2), id[u16][!][] (A1), name[str][*][C] (B1),
#   1001-100|1002-300.. (C3..C8), awards[...] (C1)
$ws.Range("C2").Value = "副本掉落"
$ws.Range("A1").Value = "id[u16][!][]"
$ws.Range("B1").Value = "name[str][*][C]"
$ws.Range("C3").Value = "1001-100|1002-300"
$ws.Range("C4").Value = "1001-100|1002-400"
$ws.Range("C5").Value = "1001-100|1002-500"
$ws.Range("C6").Value = "1001-100|1002-600"
$ws.Range("C7").Value = "1001-100|1002-700"
$ws.Range("C8").Value = "1001-100|1002-800"
$ws.Range("C1").Value = 'awards[u16-u32|][*][][Item.id-$|]'

# --- selection --------------------------------------------------------
$ws.Range("C6").Select()

# --- column widths (best-fit re-measurement after the content change) -
$ws.Columns.Item(1).ColumnWidth = 12.285714285714286
$ws.Columns.Item(2).ColumnWidth = 14.714285714285714
$ws.Columns.Item(3).ColumnWidth = 30.714285714285715

# --- page setup ---------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
